# ---------------------------------------------------------------------------
# Applies the "Staticke promenjive, formiranje SQL baze" edit:
#   - splits the trailing " subbrancha" run of paragraph 1 into a plain
#     " " run plus a spell-checked "subbrancha" run (wrapped in
#     proofErr spellStart/spellEnd, matching how Word marks it after a
#     respell)
#   - inserts a new, empty paragraph
#   - inserts a new paragraph "Instalacija sqllite3" (each word individually
#     wrapped in proofErr spellStart/spellEnd, as Word does for words it
#     doesn't recognise)
#   - inserts a final paragraph containing the sqlite tutorial URL, carrying
#     over the _GoBack bookmark that used to sit at the end of paragraph 1
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Locate the paragraph that needs to be split/extended. We search for its
# distinctive trailing word so the script keeps working even if unrelated
# content precedes it.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("subbrancha", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "edit.ps1: could not locate paragraph containing 'subbrancha'"
}

$targetPara = $searchRange.Paragraphs(1)

# Sanity-check we found the expected paragraph before rewriting anything.
$originalText = $targetPara.Range.Text
if ($originalText -notmatch "subbrancha\s*$") {
    throw "edit.ps1: unexpected paragraph text [$originalText]"
}

# The whole-document range (start of body text through the very end,
# including the paragraph mark / bookmark that close paragraph 1). Rewriting
# this complete span in one shot via InsertXML lets us precisely control
# run splitting and proofErr placement without leaving stray duplicate
# markers behind (which happens if only a partial sub-range is replaced).
$wholeDoc = $d.Range($d.Content.Start, $d.Content.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newBodyXml = (
    '<w:p ' + $ns + ' w:rsidR="00B43E45" w:rsidRDefault="00680C6F">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Inicijalizaija</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>OOPhP</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>subbrancha</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p ' + $ns + '/>' +
    '<w:p ' + $ns + '>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Instalacija</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>sqllite3</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p ' + $ns + '>' +
        '<w:r><w:t>https://www.sqlitetutorial.net/download-install-sqlite/</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/edit.xml" pkg:contentType="application/xml"><pkg:xmlData>' +
    $newBodyXml +
    '</pkg:xmlData></pkg:part></pkg:package>'

$wholeDoc.InsertXML($packageXml)
